# Add Loader for Course Table
# Populates column D (Teacher Name) for the course rows (17-58) that were
# previously left blank, using the same "testN" placeholder values the
# loader wrote, and updates the active view's top-left visible cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D17").Value = "test1"
$ws.Range("D18").Value = "test1"
$ws.Range("D19").Value = "test2"
$ws.Range("D20").Value = "test2"
$ws.Range("D21").Value = "test3"
$ws.Range("D22").Value = "test3"
$ws.Range("D23").Value = "test4"
$ws.Range("D24").Value = "test4"
$ws.Range("D25").Value = "test4"
$ws.Range("D26").Value = "test4"
$ws.Range("D27").Value = "test5"
$ws.Range("D28").Value = "test5"
$ws.Range("D29").Value = "test5"
$ws.Range("D30").Value = "test5"
$ws.Range("D31").Value = "test6"
$ws.Range("D32").Value = "test6"
$ws.Range("D33").Value = "test7"
$ws.Range("D34").Value = "test7"
$ws.Range("D35").Value = "test8"
$ws.Range("D36").Value = "test9"
$ws.Range("D37").Value = "test10"
$ws.Range("D38").Value = "test11"
$ws.Range("D39").Value = "test12"
$ws.Range("D40").Value = "test12"
$ws.Range("D41").Value = "test13"
$ws.Range("D42").Value = "test14"
$ws.Range("D43").Value = "test15"
$ws.Range("D44").Value = "test16"
$ws.Range("D45").Value = "test17"
$ws.Range("D46").Value = "test1"
$ws.Range("D47").Value = "test18"
$ws.Range("D48").Value = "test10"
$ws.Range("D49").Value = "test13"
$ws.Range("D50").Value = "test12"
$ws.Range("D51").Value = "test19"
$ws.Range("D52").Value = "test6"
$ws.Range("D53").Value = "test20"
$ws.Range("D54").Value = "test21"
$ws.Range("D55").Value = "test22"
$ws.Range("D56").Value = "test14"
$ws.Range("D57").Value = "test13"
$ws.Range("D58").Value = "test18"

# Scroll the view so row 44 is at the top-left of the visible window, while
# keeping the original selection on D59.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 44
$win.ScrollColumn = 1
$ws.Range("D59").Select() | Out-Null
